$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of notes data (row 5)
$ws.Range("A5").Value = 238
$ws.Range("B5").Value = "Product of Array Except Self O(n)"
$ws.Range("C5").Value = "Postfix and prefix"

# Update last active selection on the sheet to match the author's final click
$ws.Range("E21").Select()
